# Implement "Gyroscopic Torque" changes on the Thrust-bench test sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ------------------------------------------------------------------
# 1. Header row updates
# ------------------------------------------------------------------
# New headers for the gyroscopic-torque columns (added first so the
# shared-string table grows in the same order as the source edit)
$ws.Range("N1").Value = "K_q"
$ws.Range("O1").Value = "halved rpm K_q"
# K1 was "corrected rpm Kf" -> now "halved rpm Kf"
$ws.Range("K1").Value = "halved rpm Kf"

# ------------------------------------------------------------------
# 2. Update K_f formulas (J, K columns) to use D^4 instead of D
# ------------------------------------------------------------------
$ws.Range("J3:J12").Formula = "=I3/(C3^2*`$B`$20*`$C`$20^4)"
$ws.Range("K3:K12").Formula = "=I3*4/(C3^2*`$B`$20*`$C`$20^4)"

# ------------------------------------------------------------------
# 3. New K_q columns (N, O) - gyroscopic torque coefficient
# ------------------------------------------------------------------
$ws.Range("N2:N12").Formula = "=M2/(C2^2*`$B`$20*`$C`$20^5) * 1/(2*PI())"
# O2 was entered on its own (not part of the O3:O12 fill), matching the
# source workbook where row 2 (all-zero / #DIV0 row) was filled separately.
$ws.Range("O2").Formula = "=M2*4/(C2^2*`$B`$20*`$C`$20^5) * 1/(2*PI())"
$ws.Range("O3:O12").Formula = "=M3*4/(C3^2*`$B`$20*`$C`$20^5) * 1/(2*PI())"

# ------------------------------------------------------------------
# 4. Row 13 summary: SUM -> AVERAGE, add N13/O13 averages
# ------------------------------------------------------------------
$ws.Range("J13").Formula = "=AVERAGE(J4:J12)"
$ws.Range("K13").Formula = "=AVERAGE(K4:K12)"
$ws.Range("N13").Formula = "=AVERAGE(N4:N12)"
$ws.Range("O13").Formula = "=AVERAGE(O4:O12)"

# Keep row 13 consistently bold/filled across the now-wider used range
$ws.Range("L13:O13").Font.Bold = $true

# ------------------------------------------------------------------
# 5. New footnote / source link
# ------------------------------------------------------------------
$ws.Range("E16").Value = "https://m-selig.ae.illinois.edu/props/propDB.html"

# ------------------------------------------------------------------
# 6. Column widths (best effort)
# ------------------------------------------------------------------
$ws.Columns.Item(11).ColumnWidth = $ws.Columns.Item(10).ColumnWidth
$ws.Columns.Item(15).ColumnWidth = 13.1666666666667

# ------------------------------------------------------------------
# 7. Selection
# ------------------------------------------------------------------
$ws.Range("R19").Select()
